# Apply weekly data refresh to "Fruta, Vega Monumental Concepción - Tuna" sheet.
# Rows 2-9 get their D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de comercializacion)
# and S (Precio $/Kg) values updated per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [double]$D,
        [string]$L,
        [double]$M,
        [double]$N,
        [double]$O,
        [double]$P,
        [string]$Q,
        [double]$S
    )

    $ws.Range("D$Row").Value = $D
    $ws.Range("L$Row").Value = $L
    $ws.Range("M$Row").Value = $M
    $ws.Range("N$Row").Value = $N
    $ws.Range("O$Row").Value = $O
    $ws.Range("P$Row").Value = $P
    $ws.Range("Q$Row").Value = $Q
    $ws.Range("S$Row").Value = $S
}

Set-Row 2 44316 "Primera" 50  20000 20000 20000 "$/caja 18 kilos" 1111
Set-Row 3 44687 "Primera" 100 18000 19000 18500 "$/caja 18 kilos" 1028
Set-Row 4 44280 "Primera" 100 14000 15000 14500 "$/caja 18 kilos" 806
Set-Row 5 44280 "Segunda" 50  12000 12000 12000 "$/caja 18 kilos" 667
Set-Row 6 44699 "Primera" 100 20000 22000 21000 "$/caja 18 kilos" 1167
Set-Row 7 44699 "Segunda" 50  18000 18000 18000 "$/caja 18 kilos" 1000
Set-Row 8 44516 "Primera" 100 33000 34000 33500 "$/caja 18 kilos" 1861
Set-Row 9 44819 "Primera" 100 25000 26000 25500 "$/caja 18 kilos granel" 1417

$wb.Save()
